$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 155, shifting rows 155-246 down to 156-247.
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new data record.
$ws.Range("A155").Value = 6
$ws.Range("B155").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C155").Value = "Metropolitana"
$ws.Range("D155").Value = 44460
$ws.Range("D155").NumberFormat = $ws.Range("D156").NumberFormat
$ws.Range("E155").Value = 13
$ws.Range("F155").Value = 100112052
$ws.Range("G155").Value = "Albahaca"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 130
$ws.Range("K155").Value = 5000
$ws.Range("L155").Value = 6000
$ws.Range("M155").Value = 5538
$ws.Range("N155").Value = "$/paquete"
$ws.Range("O155").Value = "Región de Arica y Parinacota"
$ws.Range("P155").Value = 5538
$ws.Range("Q155").Value = 1
$ws.Range("R155").Value = "Hortaliza"
